# Minor edits to a couple of PowerPoint slides:
#   - Slide 3:  merge "(but "/"string literals "/"are allowed)" into a single run
#   - Slide 28: "has a method getLiteralIntValue()" -> "has a computed property intValue"
#   - Slide 29: "getLiteralIntValue()}\")" -> "intValue}\")" (two occurrences)

$p = $ppt.ActivePresentation

function Replace-FirstOccurrence {
    param($TextRange, [string]$Search, [string]$Replacement)

    $current = $TextRange.Text
    $idx = $current.IndexOf($Search)
    if ($idx -lt 0) {
        throw "Replace-FirstOccurrence: substring not found: [$Search]"
    }
    $run = $TextRange.Characters($idx + 1, $Search.Length)
    $run.Text = $Replacement
}

# ---------------------------------------------------------------------------
# Slide 28: "Class ConstValue has a method getLiteralIntValue() that returns..."
#        -> "Class ConstValue has a computed property intValue that returns..."
# ---------------------------------------------------------------------------
$slide28 = $p.Slides.Item(28)
$shape28 = $slide28.Shapes.Item(2)
$tr28 = $shape28.TextFrame.TextRange

Replace-FirstOccurrence $tr28 " has a method " " has a computed property "
Replace-FirstOccurrence $tr28 "getLiteralIntValue" "intValue"

# drop the now-orphaned "()" that followed the old method name
$afterRename = $tr28.Text
$parenIdx = $afterRename.IndexOf("intValue") + "intValue".Length
$parenRun = $tr28.Characters($parenIdx + 1, 2)
if ($parenRun.Text -eq "()") {
    $parenRun.Text = ""
}

# ---------------------------------------------------------------------------
# Slide 29: emit("LDCINT ${getLiteralIntValue()}")  -> emit("LDCINT ${intValue}")
#           emit("LDCB ${getLiteralIntValue()}")    -> emit("LDCB ${intValue}")
# ---------------------------------------------------------------------------
$slide29 = $p.Slides.Item(29)
$shape29 = $slide29.Shapes.Item(2)
$tr29 = $shape29.TextFrame.TextRange

Replace-FirstOccurrence $tr29 "getLiteralIntValue()" "intValue"
Replace-FirstOccurrence $tr29 "getLiteralIntValue()" "intValue"

# ---------------------------------------------------------------------------
# Slide 3: "(but " + "string literals " + "are allowed)" -> one merged run
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$shape3 = $slide3.Shapes.Item(2)
$tr3 = $shape3.TextFrame.TextRange

Replace-FirstOccurrence $tr3 "(but " ""
Replace-FirstOccurrence $tr3 "are allowed)" ""
Replace-FirstOccurrence $tr3 "string literals " "(but string literals are allowed)"
